# Revert responsive design implementation
# Re-adds sensor data rows that were previously trimmed from each of the
# four worksheets (ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER,
# ROW02-MID-LIFTER), plus fixes the A29 timestamp cell on ROW35-MID-LIFTER
# so it is a real date value (matching the other rows) instead of plain text.

$wb = $excel.ActiveWorkbook

# Big constant reused by every new/changed row (column G, "ID_DEC").
# Written as a full decimal literal (no exponent) so the COM layer stores
# it as a plain number without inferring a scientific-notation display
# format / extra style.
$bigG = "568631262647113769549824"

# ---------------------------------------------------------------------
# Sheet 1: ROW35-FE-LIFTER  (dimension A1:I28 -> A1:I31)
# Append rows 29-31 (same byte pattern as the rows immediately above).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(29, 1).Value = 45729.73239443287
$ws1.Cells.Item(29, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(29, 2).Value = "0x01,0x90"
$ws1.Cells.Item(29, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(29, 4).Value = "0x01,0x90,"
$ws1.Cells.Item(29, 5).Value = "0xd"
$ws1.Cells.Item(29, 6).Value = 400
$ws1.Cells.Item(29, 7).Value = $bigG
$ws1.Cells.Item(29, 8).Value = 400
$ws1.Cells.Item(29, 9).Value = 13

$ws1.Cells.Item(30, 1).Value = 45729.73241640046
$ws1.Cells.Item(30, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(30, 2).Value = "0x01,0x90"
$ws1.Cells.Item(30, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(30, 4).Value = "0x01,0x90,"
$ws1.Cells.Item(30, 5).Value = "0xd"
$ws1.Cells.Item(30, 6).Value = 400
$ws1.Cells.Item(30, 7).Value = $bigG
$ws1.Cells.Item(30, 8).Value = 400
$ws1.Cells.Item(30, 9).Value = 13

$ws1.Cells.Item(31, 1).Value = 45729.73243972223
$ws1.Cells.Item(31, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(31, 2).Value = "0x01,0x90"
$ws1.Cells.Item(31, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(31, 4).Value = "0x01,0x90,"
$ws1.Cells.Item(31, 5).Value = "0xd"
$ws1.Cells.Item(31, 6).Value = 400
$ws1.Cells.Item(31, 7).Value = $bigG
$ws1.Cells.Item(31, 8).Value = 400
$ws1.Cells.Item(31, 9).Value = 13

# ---------------------------------------------------------------------
# Sheet 2: ROW35-MID-LIFTER  (dimension A1:I29 -> A1:I32)
# Row 29's timestamp becomes a real date value; rows 30-31 are appended
# with real date values; row 32 is appended as plain text (unconverted).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(29, 1).Value = 45729.58037369213
$ws2.Cells.Item(29, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws2.Cells.Item(30, 1).Value = 45729.58039555555
$ws2.Cells.Item(30, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(30, 2).Value = "0x01,0x90"
$ws2.Cells.Item(30, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(30, 4).Value = "0x01,0x86,"
$ws2.Cells.Item(30, 5).Value = "0x4"
$ws2.Cells.Item(30, 6).Value = 400
$ws2.Cells.Item(30, 7).Value = $bigG
$ws2.Cells.Item(30, 8).Value = 390
$ws2.Cells.Item(30, 9).Value = 4

$ws2.Cells.Item(31, 1).Value = 45729.58041870371
$ws2.Cells.Item(31, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(31, 2).Value = "0x01,0x90"
$ws2.Cells.Item(31, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(31, 4).Value = "0x01,0x86,"
$ws2.Cells.Item(31, 5).Value = "0x4"
$ws2.Cells.Item(31, 6).Value = 400
$ws2.Cells.Item(31, 7).Value = $bigG
$ws2.Cells.Item(31, 8).Value = 390
$ws2.Cells.Item(31, 9).Value = 4

$ws2.Cells.Item(32, 1).Value = "2025-03-14 01:55:48"
$ws2.Cells.Item(32, 2).Value = "0x01,0x90"
$ws2.Cells.Item(32, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(32, 4).Value = "0x01,0x86,"
$ws2.Cells.Item(32, 5).Value = "0x4"
$ws2.Cells.Item(32, 6).Value = 400
$ws2.Cells.Item(32, 7).Value = $bigG
$ws2.Cells.Item(32, 8).Value = 390
$ws2.Cells.Item(32, 9).Value = 4

# ---------------------------------------------------------------------
# Sheet 3: ROW02-FE-LIFTER  (dimension A1:I31 -> A1:I32)
# Append row 32 as plain text timestamp (unconverted), like the rows
# preceding it on this sheet.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(32, 1).Value = "2025-03-13 18:26:04"
$ws3.Cells.Item(32, 2).Value = "0x01,0x90"
$ws3.Cells.Item(32, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(32, 4).Value = "0x01,0x90,"
$ws3.Cells.Item(32, 5).Value = "0x3"
$ws3.Cells.Item(32, 6).Value = 400
$ws3.Cells.Item(32, 7).Value = $bigG
$ws3.Cells.Item(32, 8).Value = 400
$ws3.Cells.Item(32, 9).Value = 3

# ---------------------------------------------------------------------
# Sheet 4: ROW02-MID-LIFTER  (dimension A1:I31 -> A1:I32)
# Append row 32 as plain text timestamp (unconverted).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(32, 1).Value = "2025-03-13 18:26:04"
$ws4.Cells.Item(32, 2).Value = "0x01,0x90"
$ws4.Cells.Item(32, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws4.Cells.Item(32, 4).Value = "0x01,0x90,"
$ws4.Cells.Item(32, 5).Value = "0x3"
$ws4.Cells.Item(32, 6).Value = 400
$ws4.Cells.Item(32, 7).Value = $bigG
$ws4.Cells.Item(32, 8).Value = 400
$ws4.Cells.Item(32, 9).Value = 3

Write-Host "Applied revert of responsive design implementation across all 4 sheets."
